$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44575
$ws.Cells.Item(2, 12).Value = 'Primera'
$ws.Cells.Item(2, 13).Value = 250
$ws.Cells.Item(2, 14).Value = 6000
$ws.Cells.Item(2, 15).Value = 6000
$ws.Cells.Item(2, 16).Value = 6000
$ws.Cells.Item(2, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(2, 19).Value = 3000

$ws.Cells.Item(3, 4).Value = 44215
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 750
$ws.Cells.Item(3, 14).Value = 4000
$ws.Cells.Item(3, 15).Value = 4000
$ws.Cells.Item(3, 16).Value = 4000
$ws.Cells.Item(3, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(3, 19).Value = 2000

$ws.Cells.Item(4, 4).Value = 44642
$ws.Cells.Item(4, 12).Value = 'Primera'
$ws.Cells.Item(4, 13).Value = 250
$ws.Cells.Item(4, 14).Value = 6000
$ws.Cells.Item(4, 15).Value = 6000
$ws.Cells.Item(4, 16).Value = 6000
$ws.Cells.Item(4, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(4, 19).Value = 3000

$ws.Cells.Item(5, 4).Value = 44239
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 350
$ws.Cells.Item(5, 14).Value = 3500
$ws.Cells.Item(5, 15).Value = 4000
$ws.Cells.Item(5, 16).Value = 3750
$ws.Cells.Item(5, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(5, 19).Value = 1875

$ws.Cells.Item(6, 4).Value = 44539
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 500
$ws.Cells.Item(6, 14).Value = 5000
$ws.Cells.Item(6, 15).Value = 5000
$ws.Cells.Item(6, 16).Value = 5000
$ws.Cells.Item(6, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(6, 19).Value = 2500

$ws.Cells.Item(7, 4).Value = 44552
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 250
$ws.Cells.Item(7, 14).Value = 6000
$ws.Cells.Item(7, 15).Value = 6000
$ws.Cells.Item(7, 16).Value = 6000
$ws.Cells.Item(7, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(7, 19).Value = 3000

$ws.Cells.Item(8, 4).Value = 44552
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 500
$ws.Cells.Item(8, 14).Value = 3600
$ws.Cells.Item(8, 15).Value = 3600
$ws.Cells.Item(8, 16).Value = 3600
$ws.Cells.Item(8, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(8, 19).Value = 1800

$ws.Cells.Item(9, 4).Value = 44167
$ws.Cells.Item(9, 12).Value = 'Primera'
$ws.Cells.Item(9, 13).Value = 250
$ws.Cells.Item(9, 14).Value = 8000
$ws.Cells.Item(9, 15).Value = 8000
$ws.Cells.Item(9, 16).Value = 8000
$ws.Cells.Item(9, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(9, 19).Value = 4000

$ws.Cells.Item(10, 4).Value = 44175
$ws.Cells.Item(10, 12).Value = 'Primera'
$ws.Cells.Item(10, 13).Value = 250
$ws.Cells.Item(10, 14).Value = 4000
$ws.Cells.Item(10, 15).Value = 4000
$ws.Cells.Item(10, 16).Value = 4000
$ws.Cells.Item(10, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(10, 19).Value = 2000

$ws.Cells.Item(11, 4).Value = 44553
$ws.Cells.Item(11, 12).Value = 'Primera'
$ws.Cells.Item(11, 13).Value = 250
$ws.Cells.Item(11, 14).Value = 6000
$ws.Cells.Item(11, 15).Value = 6000
$ws.Cells.Item(11, 16).Value = 6000
$ws.Cells.Item(11, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(11, 19).Value = 3000

$ws.Cells.Item(12, 4).Value = 44250
$ws.Cells.Item(12, 12).Value = 'Primera'
$ws.Cells.Item(12, 13).Value = 100
$ws.Cells.Item(12, 14).Value = 4000
$ws.Cells.Item(12, 15).Value = 4000
$ws.Cells.Item(12, 16).Value = 4000
$ws.Cells.Item(12, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(12, 19).Value = 2000

$ws.Cells.Item(13, 4).Value = 44582
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 13).Value = 200
$ws.Cells.Item(13, 14).Value = 6000
$ws.Cells.Item(13, 15).Value = 6000
$ws.Cells.Item(13, 16).Value = 6000
$ws.Cells.Item(13, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(13, 19).Value = 3000

$ws.Cells.Item(14, 4).Value = 44223
$ws.Cells.Item(14, 12).Value = 'Primera'
$ws.Cells.Item(14, 13).Value = 200
$ws.Cells.Item(14, 14).Value = 4000
$ws.Cells.Item(14, 15).Value = 4000
$ws.Cells.Item(14, 16).Value = 4000
$ws.Cells.Item(14, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(14, 19).Value = 2000

$ws.Cells.Item(15, 4).Value = 44253
$ws.Cells.Item(15, 12).Value = 'Primera'
$ws.Cells.Item(15, 13).Value = 25
$ws.Cells.Item(15, 14).Value = 4000
$ws.Cells.Item(15, 15).Value = 4000
$ws.Cells.Item(15, 16).Value = 4000
$ws.Cells.Item(15, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(15, 19).Value = 2000

$ws.Cells.Item(16, 4).Value = 44251
$ws.Cells.Item(16, 12).Value = 'Primera'
$ws.Cells.Item(16, 13).Value = 125
$ws.Cells.Item(16, 14).Value = 4000
$ws.Cells.Item(16, 15).Value = 4000
$ws.Cells.Item(16, 16).Value = 4000
$ws.Cells.Item(16, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(16, 19).Value = 2000

$ws.Cells.Item(17, 4).Value = 44615
$ws.Cells.Item(17, 12).Value = 'Primera'
$ws.Cells.Item(17, 13).Value = 50
$ws.Cells.Item(17, 14).Value = 6000
$ws.Cells.Item(17, 15).Value = 6000
$ws.Cells.Item(17, 16).Value = 6000
$ws.Cells.Item(17, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(17, 19).Value = 3000

$ws.Cells.Item(18, 4).Value = 44188
$ws.Cells.Item(18, 12).Value = 'Primera'
$ws.Cells.Item(18, 13).Value = 300
$ws.Cells.Item(18, 14).Value = 4000
$ws.Cells.Item(18, 15).Value = 4000
$ws.Cells.Item(18, 16).Value = 4000
$ws.Cells.Item(18, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(18, 19).Value = 2000

$ws.Cells.Item(19, 4).Value = 44188
$ws.Cells.Item(19, 12).Value = 'Primera'
$ws.Cells.Item(19, 13).Value = 500
$ws.Cells.Item(19, 14).Value = 4000
$ws.Cells.Item(19, 15).Value = 4000
$ws.Cells.Item(19, 16).Value = 4000
$ws.Cells.Item(19, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(19, 19).Value = 2000

$ws.Cells.Item(20, 4).Value = 44225
$ws.Cells.Item(20, 12).Value = 'Primera'
$ws.Cells.Item(20, 13).Value = 150
$ws.Cells.Item(20, 14).Value = 4000
$ws.Cells.Item(20, 15).Value = 4000
$ws.Cells.Item(20, 16).Value = 4000
$ws.Cells.Item(20, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(20, 19).Value = 2000

$ws.Cells.Item(21, 4).Value = 44225
$ws.Cells.Item(21, 12).Value = 'Primera'
$ws.Cells.Item(21, 13).Value = 200
$ws.Cells.Item(21, 14).Value = 4000
$ws.Cells.Item(21, 15).Value = 4000
$ws.Cells.Item(21, 16).Value = 4000
$ws.Cells.Item(21, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(21, 19).Value = 2000

$ws.Cells.Item(22, 4).Value = 44547
$ws.Cells.Item(22, 12).Value = 'Primera'
$ws.Cells.Item(22, 13).Value = 200
$ws.Cells.Item(22, 14).Value = 5000
$ws.Cells.Item(22, 15).Value = 5000
$ws.Cells.Item(22, 16).Value = 5000
$ws.Cells.Item(22, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(22, 19).Value = 2500

$ws.Cells.Item(23, 4).Value = 44616
$ws.Cells.Item(23, 12).Value = 'Primera'
$ws.Cells.Item(23, 13).Value = 100
$ws.Cells.Item(23, 14).Value = 6000
$ws.Cells.Item(23, 15).Value = 6000
$ws.Cells.Item(23, 16).Value = 6000
$ws.Cells.Item(23, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(23, 19).Value = 3000

$ws.Cells.Item(24, 4).Value = 44176
$ws.Cells.Item(24, 12).Value = 'Primera'
$ws.Cells.Item(24, 13).Value = 100
$ws.Cells.Item(24, 14).Value = 4000
$ws.Cells.Item(24, 15).Value = 4000
$ws.Cells.Item(24, 16).Value = 4000
$ws.Cells.Item(24, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(24, 19).Value = 2000

$ws.Cells.Item(25, 4).Value = 44204
$ws.Cells.Item(25, 12).Value = 'Primera'
$ws.Cells.Item(25, 13).Value = 150
$ws.Cells.Item(25, 14).Value = 4000
$ws.Cells.Item(25, 15).Value = 4000
$ws.Cells.Item(25, 16).Value = 4000
$ws.Cells.Item(25, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(25, 19).Value = 2000

$ws.Cells.Item(26, 4).Value = 44204
$ws.Cells.Item(26, 12).Value = 'Primera'
$ws.Cells.Item(26, 13).Value = 250
$ws.Cells.Item(26, 14).Value = 4000
$ws.Cells.Item(26, 15).Value = 4000
$ws.Cells.Item(26, 16).Value = 4000
$ws.Cells.Item(26, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(26, 19).Value = 2000

$ws.Cells.Item(28, 4).Value = 44224
$ws.Cells.Item(28, 12).Value = 'Primera'
$ws.Cells.Item(28, 13).Value = 250
$ws.Cells.Item(28, 14).Value = 4000
$ws.Cells.Item(28, 15).Value = 4000
$ws.Cells.Item(28, 16).Value = 4000
$ws.Cells.Item(28, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(28, 19).Value = 2000

$ws.Cells.Item(29, 4).Value = 44224
$ws.Cells.Item(29, 12).Value = 'Primera'
$ws.Cells.Item(29, 13).Value = 300
$ws.Cells.Item(29, 14).Value = 4000
$ws.Cells.Item(29, 15).Value = 4000
$ws.Cells.Item(29, 16).Value = 4000
$ws.Cells.Item(29, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(29, 19).Value = 2000

$ws.Cells.Item(30, 4).Value = 44189
$ws.Cells.Item(30, 12).Value = 'Primera'
$ws.Cells.Item(30, 13).Value = 300
$ws.Cells.Item(30, 14).Value = 3000
$ws.Cells.Item(30, 15).Value = 3000
$ws.Cells.Item(30, 16).Value = 3000
$ws.Cells.Item(30, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(30, 19).Value = 1500

$ws.Cells.Item(31, 4).Value = 44189
$ws.Cells.Item(31, 12).Value = 'Primera'
$ws.Cells.Item(31, 13).Value = 250
$ws.Cells.Item(31, 14).Value = 3000
$ws.Cells.Item(31, 15).Value = 3000
$ws.Cells.Item(31, 16).Value = 3000
$ws.Cells.Item(31, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(31, 19).Value = 1500

$ws.Cells.Item(32, 4).Value = 44540
$ws.Cells.Item(32, 12).Value = 'Primera'
$ws.Cells.Item(32, 13).Value = 250
$ws.Cells.Item(32, 14).Value = 5000
$ws.Cells.Item(32, 15).Value = 5000
$ws.Cells.Item(32, 16).Value = 5000
$ws.Cells.Item(32, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(32, 19).Value = 2500

$ws.Cells.Item(33, 4).Value = 44581
$ws.Cells.Item(33, 12).Value = 'Primera'
$ws.Cells.Item(33, 13).Value = 400
$ws.Cells.Item(33, 14).Value = 6000
$ws.Cells.Item(33, 15).Value = 6000
$ws.Cells.Item(33, 16).Value = 6000
$ws.Cells.Item(33, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(33, 19).Value = 3000

$ws.Cells.Item(34, 4).Value = 44222
$ws.Cells.Item(34, 12).Value = 'Primera'
$ws.Cells.Item(34, 13).Value = 250
$ws.Cells.Item(34, 14).Value = 4000
$ws.Cells.Item(34, 15).Value = 4000
$ws.Cells.Item(34, 16).Value = 4000
$ws.Cells.Item(34, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(34, 19).Value = 2000

$ws.Cells.Item(35, 4).Value = 44222
$ws.Cells.Item(35, 12).Value = 'Primera'
$ws.Cells.Item(35, 13).Value = 300
$ws.Cells.Item(35, 14).Value = 4000
$ws.Cells.Item(35, 15).Value = 4000
$ws.Cells.Item(35, 16).Value = 4000
$ws.Cells.Item(35, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(35, 19).Value = 2000

$ws.Cells.Item(36, 4).Value = 44201
$ws.Cells.Item(36, 12).Value = 'Primera'
$ws.Cells.Item(36, 13).Value = 200
$ws.Cells.Item(36, 14).Value = 4000
$ws.Cells.Item(36, 15).Value = 4000
$ws.Cells.Item(36, 16).Value = 4000
$ws.Cells.Item(36, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(36, 19).Value = 2000

$ws.Cells.Item(37, 4).Value = 44193
$ws.Cells.Item(37, 12).Value = 'Primera'
$ws.Cells.Item(37, 13).Value = 200
$ws.Cells.Item(37, 14).Value = 3000
$ws.Cells.Item(37, 15).Value = 3000
$ws.Cells.Item(37, 16).Value = 3000
$ws.Cells.Item(37, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(37, 19).Value = 1500

$ws.Cells.Item(38, 4).Value = 44574
$ws.Cells.Item(38, 12).Value = 'Primera'
$ws.Cells.Item(38, 13).Value = 350
$ws.Cells.Item(38, 14).Value = 6000
$ws.Cells.Item(38, 15).Value = 6000
$ws.Cells.Item(38, 16).Value = 6000
$ws.Cells.Item(38, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(38, 19).Value = 3000

$ws.Cells.Item(39, 4).Value = 44551
$ws.Cells.Item(39, 12).Value = 'Primera'
$ws.Cells.Item(39, 13).Value = 500
$ws.Cells.Item(39, 14).Value = 6000
$ws.Cells.Item(39, 15).Value = 6000
$ws.Cells.Item(39, 16).Value = 6000
$ws.Cells.Item(39, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(39, 19).Value = 3000

$ws.Cells.Item(40, 4).Value = 44544
$ws.Cells.Item(40, 12).Value = 'Primera'
$ws.Cells.Item(40, 13).Value = 250
$ws.Cells.Item(40, 14).Value = 5000
$ws.Cells.Item(40, 15).Value = 5000
$ws.Cells.Item(40, 16).Value = 5000
$ws.Cells.Item(40, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(40, 19).Value = 2500

$ws.Cells.Item(41, 4).Value = 44587
$ws.Cells.Item(41, 12).Value = 'Primera'
$ws.Cells.Item(41, 13).Value = 250
$ws.Cells.Item(41, 14).Value = 6000
$ws.Cells.Item(41, 15).Value = 6000
$ws.Cells.Item(41, 16).Value = 6000
$ws.Cells.Item(41, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(41, 19).Value = 3000

$ws.Cells.Item(42, 4).Value = 44169
$ws.Cells.Item(42, 12).Value = 'Primera'
$ws.Cells.Item(42, 13).Value = 200
$ws.Cells.Item(42, 14).Value = 5000
$ws.Cells.Item(42, 15).Value = 5000
$ws.Cells.Item(42, 16).Value = 5000
$ws.Cells.Item(42, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(42, 19).Value = 2500

$ws.Cells.Item(43, 4).Value = 44901
$ws.Cells.Item(43, 12).Value = 'Especial'
$ws.Cells.Item(43, 13).Value = 250
$ws.Cells.Item(43, 14).Value = 6000
$ws.Cells.Item(43, 15).Value = 6000
$ws.Cells.Item(43, 16).Value = 6000
$ws.Cells.Item(43, 18).Value = 'Región del Maule'
$ws.Cells.Item(43, 19).Value = 3000

$ws.Cells.Item(44, 4).Value = 44203
$ws.Cells.Item(44, 12).Value = 'Primera'
$ws.Cells.Item(44, 13).Value = 350
$ws.Cells.Item(44, 14).Value = 4000
$ws.Cells.Item(44, 15).Value = 4000
$ws.Cells.Item(44, 16).Value = 4000
$ws.Cells.Item(44, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(44, 19).Value = 2000

$ws.Cells.Item(45, 4).Value = 44558
$ws.Cells.Item(45, 12).Value = 'Primera'
$ws.Cells.Item(45, 13).Value = 100
$ws.Cells.Item(45, 14).Value = 5000
$ws.Cells.Item(45, 15).Value = 6000
$ws.Cells.Item(45, 16).Value = 5500
$ws.Cells.Item(45, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(45, 19).Value = 2750

$ws.Cells.Item(46, 4).Value = 44210
$ws.Cells.Item(46, 12).Value = 'Primera'
$ws.Cells.Item(46, 13).Value = 400
$ws.Cells.Item(46, 14).Value = 3000
$ws.Cells.Item(46, 15).Value = 4000
$ws.Cells.Item(46, 16).Value = 3500
$ws.Cells.Item(46, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(46, 19).Value = 1750

$ws.Cells.Item(47, 4).Value = 44546
$ws.Cells.Item(47, 12).Value = 'Primera'
$ws.Cells.Item(47, 13).Value = 250
$ws.Cells.Item(47, 14).Value = 5000
$ws.Cells.Item(47, 15).Value = 5000
$ws.Cells.Item(47, 16).Value = 5000
$ws.Cells.Item(47, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(47, 19).Value = 2500

$ws.Cells.Item(48, 4).Value = 44260
$ws.Cells.Item(48, 12).Value = 'Primera'
$ws.Cells.Item(48, 13).Value = 75
$ws.Cells.Item(48, 14).Value = 4000
$ws.Cells.Item(48, 15).Value = 4000
$ws.Cells.Item(48, 16).Value = 4000
$ws.Cells.Item(48, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(48, 19).Value = 2000

$ws.Cells.Item(49, 4).Value = 44586
$ws.Cells.Item(49, 12).Value = 'Primera'
$ws.Cells.Item(49, 13).Value = 350
$ws.Cells.Item(49, 14).Value = 6000
$ws.Cells.Item(49, 15).Value = 6000
$ws.Cells.Item(49, 16).Value = 6000
$ws.Cells.Item(49, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(49, 19).Value = 3000

$ws.Cells.Item(50, 4).Value = 44568
$ws.Cells.Item(50, 12).Value = 'Primera'
$ws.Cells.Item(50, 13).Value = 250
$ws.Cells.Item(50, 14).Value = 6000
$ws.Cells.Item(50, 15).Value = 6000
$ws.Cells.Item(50, 16).Value = 6000
$ws.Cells.Item(50, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(50, 19).Value = 3000

$ws.Cells.Item(51, 4).Value = 44186
$ws.Cells.Item(51, 12).Value = 'Primera'
$ws.Cells.Item(51, 13).Value = 200
$ws.Cells.Item(51, 14).Value = 4000
$ws.Cells.Item(51, 15).Value = 4000
$ws.Cells.Item(51, 16).Value = 4000
$ws.Cells.Item(51, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(51, 19).Value = 2000

$ws.Cells.Item(52, 4).Value = 44195
$ws.Cells.Item(52, 12).Value = 'Primera'
$ws.Cells.Item(52, 13).Value = 300
$ws.Cells.Item(52, 14).Value = 3000
$ws.Cells.Item(52, 15).Value = 3000
$ws.Cells.Item(52, 16).Value = 3000
$ws.Cells.Item(52, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(52, 19).Value = 1500

$ws.Cells.Item(53, 4).Value = 44567
$ws.Cells.Item(53, 12).Value = 'Primera'
$ws.Cells.Item(53, 13).Value = 250
$ws.Cells.Item(53, 14).Value = 6000
$ws.Cells.Item(53, 15).Value = 6000
$ws.Cells.Item(53, 16).Value = 6000
$ws.Cells.Item(53, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(53, 19).Value = 3000

$ws.Cells.Item(54, 4).Value = 44567
$ws.Cells.Item(54, 12).Value = 'Segunda'
$ws.Cells.Item(54, 13).Value = 250
$ws.Cells.Item(54, 14).Value = 4000
$ws.Cells.Item(54, 15).Value = 4000
$ws.Cells.Item(54, 16).Value = 4000
$ws.Cells.Item(54, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(54, 19).Value = 2000

$ws.Cells.Item(55, 4).Value = 44187
$ws.Cells.Item(55, 12).Value = 'Primera'
$ws.Cells.Item(55, 13).Value = 100
$ws.Cells.Item(55, 14).Value = 3400
$ws.Cells.Item(55, 15).Value = 3400
$ws.Cells.Item(55, 16).Value = 3400
$ws.Cells.Item(55, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(55, 19).Value = 1700

$ws.Cells.Item(56, 4).Value = 44187
$ws.Cells.Item(56, 12).Value = 'Primera'
$ws.Cells.Item(56, 13).Value = 200
$ws.Cells.Item(56, 14).Value = 4000
$ws.Cells.Item(56, 15).Value = 4000
$ws.Cells.Item(56, 16).Value = 4000
$ws.Cells.Item(56, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(56, 19).Value = 2000

$ws.Cells.Item(57, 4).Value = 44187
$ws.Cells.Item(57, 12).Value = 'Segunda'
$ws.Cells.Item(57, 13).Value = 50
$ws.Cells.Item(57, 14).Value = 3000
$ws.Cells.Item(57, 15).Value = 3000
$ws.Cells.Item(57, 16).Value = 3000
$ws.Cells.Item(57, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(57, 19).Value = 1500

$ws.Cells.Item(58, 4).Value = 44202
$ws.Cells.Item(58, 12).Value = 'Primera'
$ws.Cells.Item(58, 13).Value = 200
$ws.Cells.Item(58, 14).Value = 4000
$ws.Cells.Item(58, 15).Value = 4000
$ws.Cells.Item(58, 16).Value = 4000
$ws.Cells.Item(58, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(58, 19).Value = 2000

$ws.Cells.Item(59, 4).Value = 44610
$ws.Cells.Item(59, 12).Value = 'Primera'
$ws.Cells.Item(59, 13).Value = 250
$ws.Cells.Item(59, 14).Value = 6000
$ws.Cells.Item(59, 15).Value = 6000
$ws.Cells.Item(59, 16).Value = 6000
$ws.Cells.Item(59, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(59, 19).Value = 3000

$ws.Cells.Item(60, 4).Value = 44572
$ws.Cells.Item(60, 12).Value = 'Primera'
$ws.Cells.Item(60, 13).Value = 250
$ws.Cells.Item(60, 14).Value = 6000
$ws.Cells.Item(60, 15).Value = 6000
$ws.Cells.Item(60, 16).Value = 6000
$ws.Cells.Item(60, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(60, 19).Value = 3000

$ws.Cells.Item(61, 4).Value = 44624
$ws.Cells.Item(61, 12).Value = 'Primera'
$ws.Cells.Item(61, 13).Value = 250
$ws.Cells.Item(61, 14).Value = 6000
$ws.Cells.Item(61, 15).Value = 6000
$ws.Cells.Item(61, 16).Value = 6000
$ws.Cells.Item(61, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(61, 19).Value = 3000

$ws.Cells.Item(62, 4).Value = 44252
$ws.Cells.Item(62, 12).Value = 'Primera'
$ws.Cells.Item(62, 13).Value = 75
$ws.Cells.Item(62, 14).Value = 4000
$ws.Cells.Item(62, 15).Value = 4000
$ws.Cells.Item(62, 16).Value = 4000
$ws.Cells.Item(62, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(62, 19).Value = 2000

$ws.Cells.Item(63, 4).Value = 44614
$ws.Cells.Item(63, 12).Value = 'Primera'
$ws.Cells.Item(63, 13).Value = 300
$ws.Cells.Item(63, 14).Value = 6000
$ws.Cells.Item(63, 15).Value = 6000
$ws.Cells.Item(63, 16).Value = 6000
$ws.Cells.Item(63, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(63, 19).Value = 3000

$ws.Cells.Item(64, 4).Value = 44209
$ws.Cells.Item(64, 12).Value = 'Primera'
$ws.Cells.Item(64, 13).Value = 170
$ws.Cells.Item(64, 14).Value = 3000
$ws.Cells.Item(64, 15).Value = 4000
$ws.Cells.Item(64, 16).Value = 3500
$ws.Cells.Item(64, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(64, 19).Value = 1750

$ws.Cells.Item(65, 4).Value = 44554
$ws.Cells.Item(65, 12).Value = 'Primera'
$ws.Cells.Item(65, 13).Value = 50
$ws.Cells.Item(65, 14).Value = 6000
$ws.Cells.Item(65, 15).Value = 6000
$ws.Cells.Item(65, 16).Value = 6000
$ws.Cells.Item(65, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(65, 19).Value = 3000

$ws.Cells.Item(66, 4).Value = 44573
$ws.Cells.Item(66, 12).Value = 'Primera'
$ws.Cells.Item(66, 13).Value = 200
$ws.Cells.Item(66, 14).Value = 6000
$ws.Cells.Item(66, 15).Value = 6000
$ws.Cells.Item(66, 16).Value = 6000
$ws.Cells.Item(66, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(66, 19).Value = 3000

$ws.Cells.Item(67, 4).Value = 44609
$ws.Cells.Item(67, 12).Value = 'Primera'
$ws.Cells.Item(67, 13).Value = 400
$ws.Cells.Item(67, 14).Value = 6000
$ws.Cells.Item(67, 15).Value = 6000
$ws.Cells.Item(67, 16).Value = 6000
$ws.Cells.Item(67, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(67, 19).Value = 3000

$ws.Cells.Item(68, 4).Value = 44579
$ws.Cells.Item(68, 12).Value = 'Primera'
$ws.Cells.Item(68, 13).Value = 150
$ws.Cells.Item(68, 14).Value = 6000
$ws.Cells.Item(68, 15).Value = 6000
$ws.Cells.Item(68, 16).Value = 6000
$ws.Cells.Item(68, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(68, 19).Value = 3000

$ws.Cells.Item(69, 4).Value = 44257
$ws.Cells.Item(69, 12).Value = 'Primera'
$ws.Cells.Item(69, 13).Value = 100
$ws.Cells.Item(69, 14).Value = 4000
$ws.Cells.Item(69, 15).Value = 4000
$ws.Cells.Item(69, 16).Value = 4000
$ws.Cells.Item(69, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(69, 19).Value = 2000

$ws.Cells.Item(70, 4).Value = 44894
$ws.Cells.Item(70, 12).Value = 'Primera'
$ws.Cells.Item(70, 13).Value = 250
$ws.Cells.Item(70, 14).Value = 7000
$ws.Cells.Item(70, 15).Value = 7000
$ws.Cells.Item(70, 16).Value = 7000
$ws.Cells.Item(70, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(70, 19).Value = 3500

$ws.Cells.Item(71, 4).Value = 44221
$ws.Cells.Item(71, 12).Value = 'Primera'
$ws.Cells.Item(71, 13).Value = 150
$ws.Cells.Item(71, 14).Value = 4000
$ws.Cells.Item(71, 15).Value = 4000
$ws.Cells.Item(71, 16).Value = 4000
$ws.Cells.Item(71, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(71, 19).Value = 2000

$ws.Cells.Item(72, 4).Value = 44221
$ws.Cells.Item(72, 12).Value = 'Primera'
$ws.Cells.Item(72, 13).Value = 200
$ws.Cells.Item(72, 14).Value = 4000
$ws.Cells.Item(72, 15).Value = 4000
$ws.Cells.Item(72, 16).Value = 4000
$ws.Cells.Item(72, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(72, 19).Value = 2000

$ws.Cells.Item(73, 4).Value = 44194
$ws.Cells.Item(73, 12).Value = 'Primera'
$ws.Cells.Item(73, 13).Value = 250
$ws.Cells.Item(73, 14).Value = 4000
$ws.Cells.Item(73, 15).Value = 4000
$ws.Cells.Item(73, 16).Value = 4000
$ws.Cells.Item(73, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(73, 19).Value = 2000

$ws.Cells.Item(74, 4).Value = 44181
$ws.Cells.Item(74, 12).Value = 'Primera'
$ws.Cells.Item(74, 13).Value = 140
$ws.Cells.Item(74, 14).Value = 4000
$ws.Cells.Item(74, 15).Value = 4500
$ws.Cells.Item(74, 16).Value = 4250
$ws.Cells.Item(74, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(74, 19).Value = 2125

$ws.Cells.Item(75, 4).Value = 44566
$ws.Cells.Item(75, 12).Value = 'Primera'
$ws.Cells.Item(75, 13).Value = 250
$ws.Cells.Item(75, 14).Value = 6000
$ws.Cells.Item(75, 15).Value = 6000
$ws.Cells.Item(75, 16).Value = 6000
$ws.Cells.Item(75, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(75, 19).Value = 3000

$ws.Cells.Item(76, 4).Value = 44566
$ws.Cells.Item(76, 12).Value = 'Primera'
$ws.Cells.Item(76, 13).Value = 250
$ws.Cells.Item(76, 14).Value = 5600
$ws.Cells.Item(76, 15).Value = 5600
$ws.Cells.Item(76, 16).Value = 5600
$ws.Cells.Item(76, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(76, 19).Value = 2800

$ws.Cells.Item(77, 4).Value = 44895
$ws.Cells.Item(77, 12).Value = 'Primera'
$ws.Cells.Item(77, 13).Value = 250
$ws.Cells.Item(77, 14).Value = 7000
$ws.Cells.Item(77, 15).Value = 7000
$ws.Cells.Item(77, 16).Value = 7000
$ws.Cells.Item(77, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(77, 19).Value = 3500

$ws.Cells.Item(78, 4).Value = 44895
$ws.Cells.Item(78, 12).Value = 'Primera'
$ws.Cells.Item(78, 13).Value = 250
$ws.Cells.Item(78, 14).Value = 7000
$ws.Cells.Item(78, 15).Value = 7000
$ws.Cells.Item(78, 16).Value = 7000
$ws.Cells.Item(78, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(78, 19).Value = 3500

$ws.Cells.Item(79, 4).Value = 44211
$ws.Cells.Item(79, 12).Value = 'Primera'
$ws.Cells.Item(79, 13).Value = 200
$ws.Cells.Item(79, 14).Value = 3000
$ws.Cells.Item(79, 15).Value = 3500
$ws.Cells.Item(79, 16).Value = 3250
$ws.Cells.Item(79, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(79, 19).Value = 1625

$ws.Cells.Item(80, 4).Value = 44559
$ws.Cells.Item(80, 12).Value = 'Primera'
$ws.Cells.Item(80, 13).Value = 200
$ws.Cells.Item(80, 14).Value = 6000
$ws.Cells.Item(80, 15).Value = 6000
$ws.Cells.Item(80, 16).Value = 6000
$ws.Cells.Item(80, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(80, 19).Value = 3000

$ws.Cells.Item(81, 4).Value = 44216
$ws.Cells.Item(81, 12).Value = 'Primera'
$ws.Cells.Item(81, 13).Value = 200
$ws.Cells.Item(81, 14).Value = 4000
$ws.Cells.Item(81, 15).Value = 4000
$ws.Cells.Item(81, 16).Value = 4000
$ws.Cells.Item(81, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(81, 19).Value = 2000

$ws.Cells.Item(82, 4).Value = 44216
$ws.Cells.Item(82, 12).Value = 'Primera'
$ws.Cells.Item(82, 13).Value = 400
$ws.Cells.Item(82, 14).Value = 4000
$ws.Cells.Item(82, 15).Value = 4000
$ws.Cells.Item(82, 16).Value = 4000
$ws.Cells.Item(82, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(82, 19).Value = 2000

$ws.Cells.Item(83, 4).Value = 44217
$ws.Cells.Item(83, 12).Value = 'Primera'
$ws.Cells.Item(83, 13).Value = 250
$ws.Cells.Item(83, 14).Value = 4000
$ws.Cells.Item(83, 15).Value = 4000
$ws.Cells.Item(83, 16).Value = 4000
$ws.Cells.Item(83, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(83, 19).Value = 2000

$ws.Cells.Item(84, 4).Value = 44217
$ws.Cells.Item(84, 12).Value = 'Primera'
$ws.Cells.Item(84, 13).Value = 300
$ws.Cells.Item(84, 14).Value = 4000
$ws.Cells.Item(84, 15).Value = 4000
$ws.Cells.Item(84, 16).Value = 4000
$ws.Cells.Item(84, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(84, 19).Value = 2000

$ws.Cells.Item(85, 4).Value = 44580
$ws.Cells.Item(85, 12).Value = 'Primera'
$ws.Cells.Item(85, 13).Value = 250
$ws.Cells.Item(85, 14).Value = 6000
$ws.Cells.Item(85, 15).Value = 6000
$ws.Cells.Item(85, 16).Value = 6000
$ws.Cells.Item(85, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(85, 19).Value = 3000

$ws.Cells.Item(86, 4).Value = 44565
$ws.Cells.Item(86, 12).Value = 'Primera'
$ws.Cells.Item(86, 13).Value = 300
$ws.Cells.Item(86, 14).Value = 6000
$ws.Cells.Item(86, 15).Value = 6000
$ws.Cells.Item(86, 16).Value = 6000
$ws.Cells.Item(86, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(86, 19).Value = 3000

$ws.Cells.Item(87, 4).Value = 44565
$ws.Cells.Item(87, 12).Value = 'Primera'
$ws.Cells.Item(87, 13).Value = 250
$ws.Cells.Item(87, 14).Value = 5600
$ws.Cells.Item(87, 15).Value = 5600
$ws.Cells.Item(87, 16).Value = 5600
$ws.Cells.Item(87, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(87, 19).Value = 2800

$ws.Cells.Item(88, 4).Value = 44589
$ws.Cells.Item(88, 12).Value = 'Primera'
$ws.Cells.Item(88, 13).Value = 250
$ws.Cells.Item(88, 14).Value = 6000
$ws.Cells.Item(88, 15).Value = 6000
$ws.Cells.Item(88, 16).Value = 6000
$ws.Cells.Item(88, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(88, 19).Value = 3000

$ws.Cells.Item(89, 4).Value = 44622
$ws.Cells.Item(89, 12).Value = 'Primera'
$ws.Cells.Item(89, 13).Value = 50
$ws.Cells.Item(89, 14).Value = 6000
$ws.Cells.Item(89, 15).Value = 6000
$ws.Cells.Item(89, 16).Value = 6000
$ws.Cells.Item(89, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(89, 19).Value = 3000

$ws.Cells.Item(90, 4).Value = 44571
$ws.Cells.Item(90, 12).Value = 'Primera'
$ws.Cells.Item(90, 13).Value = 100
$ws.Cells.Item(90, 14).Value = 6000
$ws.Cells.Item(90, 15).Value = 6000
$ws.Cells.Item(90, 16).Value = 6000
$ws.Cells.Item(90, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(90, 19).Value = 3000
